$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number that Excel would
# otherwise auto-convert from text to a numeric type; force text storage
# (matching the source data's inlineStr type) via a Text number format.
# (Applied per-cell: multi-area union ranges only honour the first area
# for NumberFormat assignment in this engine.)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '70.588.67'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '3.575.39'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '585.86'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').Value = '184.99'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('D7').Value = '3.564.16'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.621'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('E10').Value = '  +16.85%  '
$ws.Range('D11').Value = '0.650'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '54.11'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('E13').Value = '  +5.61%  '
$ws.Range('D14').Value = '9.55'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '4.150.94'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').Value = '19.52'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').Value = '70.644.61'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '3.577.80'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '572.19'
$ws.Range('E19').Value = '  +15.67%  '
$ws.Range('D20').Value = '12.36'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').Value = '17.66'
$ws.Range('E23').Value = '  -9.92%  '
$ws.Range('D24').Value = '4.65'
$ws.Range('E24').Value = '  +5.41%  '
$ws.Range('D25').Value = '4.92'
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('D26').Value = '94.80'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').Value = '11.35'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').Value = '2.94'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('D29').Value = '9.12'
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('D30').Value = '32.24'
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('D31').Value = '7.28'
$ws.Range('E31').Value = '  -5.94%  '
$ws.Range('D32').Value = '12.31'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = '64.67'
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('D36').Value = '561.69'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '37.60'
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('D40').Value = '0.0₃0791'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = '3.391.34'
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '3.10'
$ws.Range('E43').Value = '  -2.69%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '3.37'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('D45').Value = '3.59'
$ws.Range('E45').Value = '  -4.34%  '
$ws.Range('D46').Value = '0.0445'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').Value = '2.96'
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = '1.46'
$ws.Range('E51').Value = '  -3.24%  '
